# Apply "last report 02-02-25" update to the Route Cost RSO workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Route")

# L3 previously held the text "28/1/2025" as a shared string; replace it
# with the actual date value for 2025-02-02 (the cell is already formatted
# as a date via its style, numFmtId 14 => m/d/yyyy).
$ws.Range("L3").Value = Get-Date -Year 2025 -Month 2 -Day 2 -Hour 0 -Minute 0 -Second 0

# Row 7 (RSO 01 / Asim Gain) unit cost reduced from 150 to 130.
$ws.Range("D7").Value = 130

# Row 9 (RSO 03 / Liton Ray) unit cost reduced from 150 to 130.
$ws.Range("D9").Value = 130

$wb.Save()
